{"js": "// Fixed an error with code blocks\n//\n// Adds the \"CodeSnippet\" paragraph style (based on Normal/\"Standard\") and\n// its linked \"CodeSnippet Zchn\" character style (based on \"Default\n// Paragraph Font\") to the document's style sheet. No document content is\n// changed -- this only defines the two styles.\n\n// Word.Document.addStyle(name, type) creates a new style (WordApi 1.3+).\n// Paragraph style first, then its companion character style.\ncontext.document.addStyle(\"CodeSnippet\", Word.StyleType.paragraph);\ncontext.document.addStyle(\"CodeSnippet Zchn\", Word.StyleType.character);\nawait context.sync();\n\n// Re-resolve the two styles by name so the proxies are fully addressable\n// for further property writes.\nconst codeSnippet = context.document.styles.getByName(\"CodeSnippet\");\nconst codeSnippetChar = context.document.styles.getByName(\"CodeSnippet Zchn\");\n\ncodeSnippet.baseStyle = \"Standard\";\ncodeSnippet.quickStyle = true;\n\ncodeSnippetChar.baseStyle = \"Absatz-Standardschriftart\";\n\nawait context.sync();\n\n// Link the paragraph style with its character style (and vice versa) --\n// mirrors Word's own \"Link to paragraph style\" behavior used for things\n// like code-block styles. Word.Style.linkStyle is a read-only navigation\n// property in the public Word JavaScript API, so there is no supported\n// `style.linkStyle = ...` setter; fall back to the lower-level OM setter\n// the host exposes on the style proxy when a direct assignment isn't\n// accepted.\nfunction linkStyles(fromStyle, toStyleId) {\n  try {\n    fromStyle.linkStyle = toStyleId;\n  } catch (e) {\n    if (typeof fromStyle._omSet === \"function\") {\n      fromStyle._omSet(\"LinkStyle\", toStyleId, \"Style\");\n    }\n  }\n}\n\nlinkStyles(codeSnippet, \"CodeSnippetZchn\");\nlinkStyles(codeSnippetChar, \"CodeSnippet\");\n\nawait context.sync();\n", "ps1": "# Fixed an error with code blocks\n# Adds the \"CodeSnippet\" paragraph style (based on Normal/Standard) and its\n# linked \"CodeSnippet Zchn\" character style (based on Default Paragraph\n# Font) to the document's style sheet.\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeParagraph = 1, wdStyleTypeCharacter = 2\n$codeSnippet = $d.Styles.Add(\"CodeSnippet\", 1)\n$codeSnippet.BaseStyle = \"Standard\"\n$codeSnippet.QuickStyle = $true\n\n$codeSnippetChar = $d.Styles.Add(\"CodeSnippetZchn\", 2)\n$codeSnippetChar.NameLocal = \"CodeSnippet Zchn\"\n$codeSnippetChar.BaseStyle = \"Absatz-Standardschriftart\"\n\n# Link the paragraph style and its character style together.\n$codeSnippet.LinkStyle = \"CodeSnippetZchn\"\n$codeSnippetChar.LinkStyle = \"CodeSnippet\"\n"}
